$d = $word.ActiveDocument

# 1. Replace the text of the paragraph that previously read
#    "These cause of these bugs are not well known..." with the new
#    "-Insert (Deletions)" bullet point.
$d.Content.Find.Execute(
    "These cause of these bugs are not well known, making it difficult to fix. None of these bugs crash the program however, so no data is lost.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-Insert (Deletions): Not implemented, if you delete characters while typing others (by first pressing inser), it will throw the entire history out of sync.",
    2) | Out-Null

# 2. The following two blank paragraphs collapse into one: the first
#    blank paragraph now carries the "These cause of some of these bugs..."
#    text, and the second (now redundant) blank paragraph is removed.
$p21 = $d.Paragraphs.Item(21)
$p22 = $d.Paragraphs.Item(22)

$p21.Range.Text = "These cause of some of these bugs are not well known, making it difficult to fix. None of these bugs crash the program however, so no data is lost."

$p22.Range.Delete()
